# Adds I0 and IF columns (I, J) to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-39: [I, J] values
$values = @(
    @(1, 5),
    @(1, 7),
    @(1, 7),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 5),
    @(1, 3),
    @(1, 6),
    @(7, 9),
    @(5, 6),
    @(7, 7),
    @(6, 6),
    @(9, 9),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(5, 6),
    @(6, 6),
    @(5, 5),
    @(7, 8),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
